# Update "想去人数" (number of people interested) values in both the
# "展览" sheet and the aggregated "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    2  = 6667
    6  = 2066
    7  = 1575
    8  = 313
    10 = 455
    12 = 5647
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
